$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 5697.143
$ws.Range("I62").Value = 1970
$ws.Range("K62").Value = 1970
$ws.Range("M62").Value = -1346

# Row 65
$ws.Range("H65").Value = 5697.143
$ws.Range("I65").Value = 1970
$ws.Range("K65").Value = 9850
$ws.Range("M65").Value = -6730

# Row 132
$ws.Range("H132").Value = 100431.63
$ws.Range("I132").Value = 116268.484
$ws.Range("J132").Value = 8050
$ws.Range("K132").Value = 348805.452
$ws.Range("L132").Value = 24150
$ws.Range("M132").Value = -346275.452
$ws.Range("N132").Value = -29210

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7396.9775
$ws.Range("I32").Value = 4462.7446
$ws.Range("K32").Value = 4462.7446
$ws.Range("M32").Value = -4175.7446

# Row 53
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 137
$ws.Range("H137").Value = 51171.6
$ws.Range("J137").Value = 51171.6
$ws.Range("L137").Value = 51171.6
$ws.Range("N137").Value = -61371.6

# Row 139
$ws.Range("H139").Value = 41172.082
$ws.Range("J139").Value = 41172.082
$ws.Range("L139").Value = 41172.082
$ws.Range("N139").Value = -51452.082

$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 5493.2
$ws.Range("I54").Value = 1366.5
$ws.Range("K54").Value = 1366.5
$ws.Range("M54").Value = -882.5

# Row 138
$ws.Range("H138").Value = 41080
$ws.Range("J138").Value = 41080
$ws.Range("L138").Value = 41080
$ws.Range("N138").Value = -51360

$ws = $wb.Worksheets.Item("CRP")
# Row 54
$ws.Range("H54").Value = 10092
$ws.Range("J54").Value = 10092
$ws.Range("L54").Value = 10092
$ws.Range("N54").Value = -11408

# Row 58
$ws.Range("H58").Value = 2048.2742
$ws.Range("I58").Value = 1788.7358
$ws.Range("J58").Value = 3576.6667
$ws.Range("K58").Value = 1788.7358
$ws.Range("L58").Value = 3576.6667
$ws.Range("M58").Value = -1585.7358
$ws.Range("N58").Value = -3982.6667

# Row 86
$ws.Range("H86").Value = 2962.5
$ws.Range("I86").Value = 2350
$ws.Range("K86").Value = 2350
$ws.Range("M86").Value = -1227

# Row 89
$ws.Range("H89").Value = 2962.5
$ws.Range("I89").Value = 2350
$ws.Range("K89").Value = 11750
$ws.Range("M89").Value = -6134

# Row 132
$ws.Range("H132").Value = 3375.9524
$ws.Range("I132").Value = 2494.7144
$ws.Range("J132").Value = 5138.4287
$ws.Range("K132").Value = 7484.1432
$ws.Range("L132").Value = 15415.2861
$ws.Range("M132").Value = -4954.1432
$ws.Range("N132").Value = -20475.2861

# Row 136
$ws.Range("H136").Value = 2048.2742
$ws.Range("I136").Value = 1788.7358
$ws.Range("J136").Value = 3576.6667
$ws.Range("K136").Value = 5366.207399999999
$ws.Range("L136").Value = 10730.0001
$ws.Range("M136").Value = -2816.207399999999
$ws.Range("N136").Value = -15830.0001

# Row 138
$ws.Range("H138").Value = 34148.89
$ws.Range("J138").Value = 34148.89
$ws.Range("L138").Value = 34148.89
$ws.Range("N138").Value = -44428.89

# Row 140
$ws.Range("H140").Value = 116617.5
$ws.Range("J140").Value = 116617.5
$ws.Range("L140").Value = 116617.5
$ws.Range("N140").Value = -126977.5

$ws = $wb.Worksheets.Item("CUL")
# Row 47
$ws.Range("H47").Value = 1750
$ws.Range("I47").Value = 500
$ws.Range("J47").Value = 2000
$ws.Range("K47").Value = 1500
$ws.Range("L47").Value = 6000
$ws.Range("M47").Value = -1069
$ws.Range("N47").Value = -6862

# Row 48
$ws.Range("H48").Value = 7120
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 7120
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 21360
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -21860

# Row 131
$ws.Range("H131").Value = 13514515
$ws.Range("J131").Value = 793.24243
$ws.Range("L131").Value = 2379.72729
$ws.Range("N131").Value = -12459.72729

# Row 132
$ws.Range("H132").Value = 2234.5667
$ws.Range("I132").Value = 880.3889
$ws.Range("J132").Value = 4265.8335
$ws.Range("K132").Value = 7923.5001
$ws.Range("L132").Value = 38392.5015
$ws.Range("M132").Value = -5393.5001
$ws.Range("N132").Value = -43452.5015

$ws = $wb.Worksheets.Item("GSM")
# Row 48
$ws.Range("H48").Value = 35000
$ws.Range("J48").Value = 35000
$ws.Range("L48").Value = 35000
$ws.Range("N48").Value = -35970

# Row 53
$ws.Range("H53").Value = 26998
$ws.Range("J53").Value = 26998
$ws.Range("L53").Value = 26998
$ws.Range("N53").Value = -28260

# Row 55
$ws.Range("H55").Value = 30000
$ws.Range("J55").Value = 30000
$ws.Range("L55").Value = 30000
$ws.Range("N55").Value = -30654

# Row 137
$ws.Range("H137").Value = 72206.44500000001
$ws.Range("J137").Value = 72206.44500000001
$ws.Range("L137").Value = 72206.44500000001
$ws.Range("N137").Value = -82406.44500000001

# Row 140
$ws.Range("H140").Value = 42278
$ws.Range("J140").Value = 42278
$ws.Range("L140").Value = 42278
$ws.Range("N140").Value = -52638

$ws = $wb.Worksheets.Item("LTW")
# Row 53
$ws.Range("H53").Value = 31525.5
$ws.Range("J53").Value = 31525.5
$ws.Range("L53").Value = 31525.5
$ws.Range("N53").Value = -32561.5

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 55
$ws.Range("H55").Value = 291.1579
$ws.Range("I55").Value = 215.66667
$ws.Range("K55").Value = 215.66667
$ws.Range("M55").Value = -42.66667000000001

# Row 136
$ws.Range("H136").Value = 3008.1777
$ws.Range("I136").Value = 1400.4062
$ws.Range("J136").Value = 6965.769
$ws.Range("K136").Value = 4201.2186
$ws.Range("L136").Value = 20897.307
$ws.Range("M136").Value = -1651.2186
$ws.Range("N136").Value = -25997.307

# Row 140
$ws.Range("H140").Value = 66662.27
$ws.Range("J140").Value = 66662.27
$ws.Range("L140").Value = 66662.27
$ws.Range("N140").Value = -77022.27

# Row 141
$ws.Range("H141").Value = 41495
$ws.Range("J141").Value = 41495
$ws.Range("L141").Value = 41495

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 3684.5908
$ws.Range("I126").Value = 2678.7693
$ws.Range("J126").Value = 5137.4443
$ws.Range("K126").Value = 8036.3079
$ws.Range("L126").Value = 15412.3329
$ws.Range("M126").Value = -5566.3079
$ws.Range("N126").Value = -20352.3329

# Row 133
$ws.Range("H133").Value = 61500
$ws.Range("J133").Value = 61500
$ws.Range("L133").Value = 61500
$ws.Range("N133").Value = -71620

# Row 139
$ws.Range("H139").Value = 41081.363
$ws.Range("J139").Value = 41101.906
$ws.Range("L139").Value = 41101.906
$ws.Range("N139").Value = -51381.906

# Row 141
$ws.Range("H141").Value = 44084.117
$ws.Range("J141").Value = 44084.117
$ws.Range("L141").Value = 44084.117
$ws.Range("N141").Value = -54444.117
